# Auto update ESO load data
# Refresh forecast/actual load rows 6-10 (Feb 9-13, 2026) with latest ESO figures
# and append the new day row 11 (Feb 14, 2026).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 6
$loadRows = @(
    @(2026,2,9,5048,4722,4574,4581,4586,4678,4956,5349,5757,5990,6053,6028,5963,5923,5849,5703,5590,5736,6165,6154,6012,5806,5604,5363),
    @(2026,2,10,5374,5082,4948,4942,4925,5041,5407,5822,6180,6328,6249,6114,5979,5898,5782,5778,5813,6026,6431,6406,6254,6044,5884,5644),
    @(2026,2,11,5343,5051,4917,4911,4894,5010,5376,5791,6149,6296,6219,6086,5953,5874,5760,5756,5791,6001,6400,6375,6223,6012,5852,5612),
    @(2026,2,12,5131,4838,4703,4698,4680,4796,5162,5577,5934,6081,6013,5897,5782,5712,5613,5609,5641,5831,6192,6167,6012,5798,5636,5392),
    @(2026,2,13,4920,4626,4490,4484,4467,4583,4948,5362,5719,5866,5808,5709,5610,5550,5466,5462,5487,5631,5904,5879,5728,5520,5361,5123),
    @(2026,2,14,4785,4533,4410,4279,4299,4395,4529,4696,5055,5278,5422,5349,5321,5231,5166,5133,5167,5380,5577,5478,5315,5119,5053,4894)
)

for ($rIdx = 0; $rIdx -lt $loadRows.Length; $rIdx++) {
    $rowNum = $startRow + $rIdx
    $values = $loadRows[$rIdx]
    for ($cIdx = 0; $cIdx -lt $values.Length; $cIdx++) {
        $ws.Cells.Item($rowNum, $cIdx + 1).Value = $values[$cIdx]
    }
}

